# Fruta / hortaliza, semanal
# Insert a new weekly record at row 314 (Feria Lagunitas de Puerto Montt - Ciboulette),
# which pushes the existing rows 314:339 down to 315:340.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 314, shifting the rest of the table down.
$ws.Range("A314").EntireRow.Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A314").Value = 4
$ws.Range("B314").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C314").Value = "Los Lagos"
$ws.Range("D314").Value = 45013
$ws.Range("E314").Value = 10
$ws.Range("F314").Value = 100112039
$ws.Range("G314").Value = "Ciboulette"
$ws.Range("H314").Value = "Sin especificar"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 240
$ws.Range("K314").Value = 3500
$ws.Range("L314").Value = 3500
$ws.Range("M314").Value = 3500
$ws.Range("N314").Value = "$/docena de atados"
$ws.Range("O314").Value = "Región Metropolitana"
$ws.Range("P314").Value = 1167
$ws.Range("Q314").Value = 3
$ws.Range("R314").Value = "Hortaliza"
